# Rotate the full row contents (columns A:AY) among rows 8-16.
# The mapping below means: the data that *ends up* in the target row (key)
# is the data that was originally in the source row (value), before any
# writes happen. We snapshot every source row first, then write them all
# back, so overlapping reads/writes never clobber data we still need.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 8
$lastRow = 16
$colRange = "A{0}:AY{0}"

# 1) Snapshot current contents of rows 8..16 (each as a 1 x 51 Value2 array)
#    plus the plain-text date cells (Y, AA) that need special handling so
#    Excel doesn't silently reinterpret them as date serials on write-back.
$snapshot = @{}
$ySnapshot = @{}
$aaSnapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $addr = ($colRange -f $r)
    $snapshot[$r] = $ws.Range($addr).Value2
    $ySnapshot[$r] = $ws.Range("Y$r").Value2
    $aaSnapshot[$r] = $ws.Range("AA$r").Value2
}

# 2) Target row -> source row mapping (content moves from source to target).
$mapping = @{
    8  = 12
    9  = 13
    10 = 14
    11 = 15
    12 = 8
    13 = 16
    14 = 9
    15 = 10
    16 = 11
}

# 3) Write each target row using the snapshot taken before any writes.
foreach ($target in ($mapping.Keys | Sort-Object)) {
    $source = $mapping[$target]
    $addr = ($colRange -f $target)
    $ws.Range($addr).Value = $snapshot[$source]
}

# 4) Columns Y and AA hold plain text dates like "2022-06-02". A plain
#    array/.Value write lets Excel auto-detect and convert such strings
#    into real date serials, which changes both the stored value and its
#    type. Force them to literal text (matching the original encoding),
#    then drop the temporary formatting so no stray styling is left behind.
$ws.Range("Y$firstRow`:Y$lastRow").NumberFormat = "@"
$ws.Range("AA$firstRow`:AA$lastRow").NumberFormat = "@"
foreach ($target in ($mapping.Keys | Sort-Object)) {
    $source = $mapping[$target]
    $ws.Range("Y$target").Value = $ySnapshot[$source]
    $ws.Range("AA$target").Value = $aaSnapshot[$source]
}
$ws.Range("Y$firstRow`:Y$lastRow").Style = "Normal"
$ws.Range("AA$firstRow`:AA$lastRow").Style = "Normal"
